$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 13, shifting existing rows 13-23 down to 14-24 ---
$ws.Rows.Item(13).Insert()

# --- Row 10: Objetivos answer (now filled with the actual objectives text) ---
$ws.Range("B10").Value = "Capacitar os alunos a calcular os parâmetros de projeto de reatores ideais, a distinguir entre um reator ideal e um real, e a compreender a influência da temperatura e pressão no projeto de reatores químicos."
$ws.Range("C10").Value = "Capacitar os alunos a calcular os parâmetros de projeto de reatores ideais, a distinguir entre um reator ideal e um real, e a compreender a influência da temperatura e pressão no projeto de reatores químicos."

# --- Row 13 (new, blank label cell; B/C hold the professor name, moved up from old row 10) ---
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "5963230 - Leandro Gonçalves de Aguiar"
$ws.Range("C13").Value = "5963230 - Leandro Gonçalves de Aguiar"
$ws.Rows.Item(13).EntireRow.AutoFit()

# --- Row 14: "Programa resumido:" now has its proper short-syllabus summary text ---
$ws.Range("B14").Value = "1. Introdução a Reatores. 2. Modelos Ideais de Reatores Químicos Isotérmicos  Reações Simples. 3. Reações Múltiplas em Reatores Ideais. 4. Efeitos Térmicos em Reatores Ideais. 5. Reatores Catalíticos Heterogêneos. 6. Reatores Não-Ideais"
$ws.Range("C14").Value = "1. Introdução a Reatores. 2. Modelos Ideais de Reatores Químicos Isotérmicos  Reações Simples. 3. Reações Múltiplas em Reatores Ideais. 4. Efeitos Térmicos em Reatores Ideais. 5. Reatores Catalíticos Heterogêneos. 6. Reatores Não-Ideais"

# --- Row 15: "Short syllabus:" label only; clear stray date that had been here ---
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: "Programa:" gets the full detailed syllabus text ---
$ws.Range("B16").Value = "1. Introdução a Reatores: Conceitos básicos`n2. Modelos Ideais de Reatores Químicos Isotérmicos  Reações Simples: `n2.1) Equações fundamentais de projeto de reatores`n2.2) Reator tanque descontínuo (BSTR)`n2.3) Reator tanque de mistura contínuo (CSTR)`n2.4) Reator tubular de fluxo pistonado (PFR)`n2.5) Comparação de desempenho de reatores CSTR e PFR`n2.6) Reatores CSTR em cascata`n2.7) Associação mista de reatores em série: CSTR e PFR`n2.8) Reatores com reciclo`n2.9) Reações auto-catalíticas`n2.10) Reatores semi-contínuos`n3. Reações Múltiplas em Reatores Ideais`n3.1) Noções gerais: otimização, rendimento e seletividade`n3.2) Reações paralelas e reações em série`n3.3) Sistemas com reações série-paralelo: reações de múltipla substituição e reações poliméricas`n3.4) Problemas simples de otimização`n4. Efeitos Térmicos em Reatores Ideais`n4.1) Equação do balanço de energia`n4.2) Balanço de energia aplicado ao BSTR`n4.3) Balanço de energia aplicado ao CSTR`n4.4) Balanço de energia aplicado ao PFR`n5. Reatores Catalíticos Heterogêneos`n5.1) Introdução`n5.2) Efeito dos processos físicos sobre a taxa de reação`n5.2.1  Fenômenos interfases`n5.2.2  Fenômenos intrapartícula`n5.2.3  Difusão e reação em catalisadores porosos`n5.3) Cálculo de reatores de leito fixo`n5.4) Reatores trifásicos`n6. Reatores Não-Ideais`n6.1) A distribuição dos tempos de residência`n6.2) Modelos dos tanques contínuos em série`n6.3) Modelo da dispersão axial"
$ws.Range("C16").Value = "1. Introdução a Reatores: Conceitos básicos`n2. Modelos Ideais de Reatores Químicos Isotérmicos  Reações Simples: `n2.1) Equações fundamentais de projeto de reatores`n2.2) Reator tanque descontínuo (BSTR)`n2.3) Reator tanque de mistura contínuo (CSTR)`n2.4) Reator tubular de fluxo pistonado (PFR)`n2.5) Comparação de desempenho de reatores CSTR e PFR`n2.6) Reatores CSTR em cascata`n2.7) Associação mista de reatores em série: CSTR e PFR`n2.8) Reatores com reciclo`n2.9) Reações auto-catalíticas`n2.10) Reatores semi-contínuos`n3. Reações Múltiplas em Reatores Ideais`n3.1) Noções gerais: otimização, rendimento e seletividade`n3.2) Reações paralelas e reações em série`n3.3) Sistemas com reações série-paralelo: reações de múltipla substituição e reações poliméricas`n3.4) Problemas simples de otimização`n4. Efeitos Térmicos em Reatores Ideais`n4.1) Equação do balanço de energia`n4.2) Balanço de energia aplicado ao BSTR`n4.3) Balanço de energia aplicado ao CSTR`n4.4) Balanço de energia aplicado ao PFR`n5. Reatores Catalíticos Heterogêneos`n5.1) Introdução`n5.2) Efeito dos processos físicos sobre a taxa de reação`n5.2.1  Fenômenos interfases`n5.2.2  Fenômenos intrapartícula`n5.2.3  Difusão e reação em catalisadores porosos`n5.3) Cálculo de reatores de leito fixo`n5.4) Reatores trifásicos`n6. Reatores Não-Ideais`n6.1) A distribuição dos tempos de residência`n6.2) Modelos dos tanques contínuos em série`n6.3) Modelo da dispersão axial"

# --- Row 17: "Syllabus:" label only ---
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: "Avaliação:" label only; clear stray professor name that had been here ---
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Rows.Item(18).EntireRow.AutoFit()

# --- Row 19: "Método:" now holds the exam-method text ---
$ws.Range("B19").Value = "Duas provas escritas e eventual apresentação de trabalho."
$ws.Range("C19").Value = "Duas provas escritas e eventual apresentação de trabalho."

# --- Row 20: "Critério:" now holds the grading-criteria text ---
$ws.Range("B20").Value = "Nota(N) = 50% Prova P1 + 50% Prova P2. Os pesos poderão ser redefinidos caso seja incorporada nota de trabalho."
$ws.Range("C20").Value = "Nota(N) = 50% Prova P1 + 50% Prova P2. Os pesos poderão ser redefinidos caso seja incorporada nota de trabalho."

# --- Row 21: "Norma de recuperação:" now holds the make-up exam rule text ---
$ws.Range("B21").Value = "Média Final = (N + Prova Recuperação)/2"
$ws.Range("C21").Value = "Média Final = (N + Prova Recuperação)/2"
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22: "Bibliografia:" now holds the bibliography text ---
$ws.Range("B22").Value = "FOGLER, H. S. Elementos de Engenharia das Reações Químicas. 3. ed. Rio de Janeiro: LTC Editora, 2002.`nLEVENSPIEL, O. Chemical Reaction Engineering. 3. ed. New York: John Wiley & Sons, 1998.`nHILL, C.G. An Introduction to Chemical Engineering Kinetics and Reactor Design. New York: John Wiley&Sons, 1977.`nSMITH, J.M. Chemical Engineering Kinetics. 3rd. ed. New York :  McGraw-Hill, 1981.`nDENBIGH, K.; TURNER, R. Introduction to Chemical Reaction Design. Cambridge: Cambridge University Press, 1970.`nFROMENT, G.F.; BISCHOFF, K.B. Chemical Reactor Analysis And Design. 2nd ed.  New York: John Wiley & Sons, 1990."
$ws.Range("C22").Value = "FOGLER, H. S. Elementos de Engenharia das Reações Químicas. 3. ed. Rio de Janeiro: LTC Editora, 2002.`nLEVENSPIEL, O. Chemical Reaction Engineering. 3. ed. New York: John Wiley & Sons, 1998.`nHILL, C.G. An Introduction to Chemical Engineering Kinetics and Reactor Design. New York: John Wiley&Sons, 1977.`nSMITH, J.M. Chemical Engineering Kinetics. 3rd. ed. New York :  McGraw-Hill, 1981.`nDENBIGH, K.; TURNER, R. Introduction to Chemical Reaction Design. Cambridge: Cambridge University Press, 1970.`nFROMENT, G.F.; BISCHOFF, K.B. Chemical Reactor Analysis And Design. 2nd ed.  New York: John Wiley & Sons, 1990."
$ws.Rows.Item(22).RowHeight = 120

# --- Row 23: "Requisitos:" label only; clear the prerequisite text that had been here ---
$ws.Range("B23").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Rows.Item(23).EntireRow.AutoFit()

# --- Row 24 (new last row): prerequisite text, pushed down from old row 23 ---
$ws.Range("B24").Value = "LOQ4003 -  Cinética Química Aplicada  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOQ4003 -  Cinética Química Aplicada  (Requisito fraco)`n"
$ws.Rows.Item(24).RowHeight = 30

# --- Column layout: column A width no longer spans into column B ---
$ws.Columns.Item(1).ColumnWidth = 30.7109375
$ws.Columns.Item(2).ColumnWidth = 60.7109375
$ws.Columns.Item(3).ColumnWidth = 60.7109375
